$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.501.82"
$ws.Range("E2").Value = "  +2.48%  "

$ws.Range("D3").Value = "1.678.95"

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.48"
$ws.Range("E5").Value = "  +3.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5322"
$ws.Range("E6").Value = "  +2.31%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2678"
$ws.Range("E8").Value = "  +4.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06394"
$ws.Range("E9").Value = "  +1.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.63"
$ws.Range("E10").Value = "  +5.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07803"
$ws.Range("E11").Value = "  +3.39%  "

$ws.Range("D12").Value = "1.691.08"
$ws.Range("E12").Value = "  +4.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.500"
$ws.Range("E13").Value = "  +2.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5569"
$ws.Range("E14").Value = "  +1.32%  "

$ws.Range("D15").Value = "0.0₅8338"
$ws.Range("E15").Value = "  +4.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.68"
$ws.Range("E16").Value = "  +1.97%  "

$ws.Range("D17").Value = "26.530.88"
$ws.Range("E17").Value = "  +2.65%  "

$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.765"
$ws.Range("E19").Value = "  +2.47%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.69"
$ws.Range("E20").Value = "  +5.63%  "

$ws.Range("E21").Value = "  +2.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.342"
$ws.Range("E22").Value = "  +4.18%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.49"
$ws.Range("E24").Value = "  -0.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1283"
$ws.Range("E25").Value = "  +5.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.430"
$ws.Range("E26").Value = "  +0.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.35"
$ws.Range("E27").Value = "  +4.96%  "

$ws.Range("E28").Value = "  +5.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06166"
$ws.Range("E29").Value = "  +5.17%  "

$ws.Range("E30").Value = "  +2.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.609"
$ws.Range("E31").Value = "  +6.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.461"
$ws.Range("E32").Value = "  +2.85%  "

$ws.Range("E33").Value = "  +4.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.008"
$ws.Range("E34").Value = "  +3.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.424"
$ws.Range("E35").Value = "  +1.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.786"
$ws.Range("E36").Value = "  +1.86%  "

$ws.Range("E37").Value = "  -0.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01638"
$ws.Range("E38").Value = "  +2.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.036"
$ws.Range("E39").Value = "  +6.85%  "

$ws.Range("D40").Value = "1.077.67"
$ws.Range("E40").Value = "  +4.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8595"
$ws.Range("E41").Value = "  +1.68%  "

$ws.Range("E42").Value = "  -0.29%  "

$ws.Range("E43").Value = "  +0.41%  "

$ws.Range("D44").Value = "1.824.01"
$ws.Range("E44").Value = "  +3.01%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₈109"
$ws.Range("E45").Value = "  +0.55%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.00"
$ws.Range("E46").Value = "  +4.16%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.168"
$ws.Range("E47").Value = "  +2.18%  "

$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  +0.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05210"
$ws.Range("E49").Value = "  +1.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.032"
$ws.Range("E50").Value = "  +3.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4236"
$ws.Range("E51").Value = "  +0.58%  "
